$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": insert a new date column before DT (shifts DT:EX
# right to DU:EY) and populate it with the "19-nov" header plus a "-"
# placeholder for every data row (2-25), matching the rest of the
# not-yet-started dates already present on the sheet.
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("DT1").EntireColumn.Insert()

$wsPrix.Range("DT1").Value = "19-nov"
$wsPrix.Range("DT2:DT25").Value = "-"

# ---------------------------------------------------------------------
# Sheet "Gaz": append a new row with the 2025-11-17 price point.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A153").NumberFormat = "@"
$wsGaz.Range("A153").Value = "2025-11-17"
$wsGaz.Range("A153").ClearFormats()
$wsGaz.Range("B153").Value = 30.395

# ---------------------------------------------------------------------
# Sheet "CO2": append a new row with the 2025-11-17 price point.
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A153").NumberFormat = "@"
$wsCo2.Range("A153").Value = "2025-11-17"
$wsCo2.Range("A153").ClearFormats()
$wsCo2.Range("B153").Value = 79.68000000000001
